$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C28").Value = 3979.85302734375
$ws.Range("E28").Value = 0.00000000054060328436023042
$ws.Range("C29").Value = 5969.779296875
$ws.Range("E29").Value = 0.00000000071604633333777201
$ws.Range("C30").Value = 11939.8681640625
$ws.Range("E30").Value = 0.0000000036590239727019025
$ws.Range("C31").Value = 23879.921875
$ws.Range("E31").Value = 0.0000000076143660265870494
$ws.Range("C32").Value = 45768.3125
$ws.Range("E32").Value = 0.000000019715924892693693
$ws.Range("C33").Value = 77610.2578125
$ws.Range("E33").Value = 0.000000024852688795817812
$ws.Range("C34").Value = 117786.9765625
$ws.Range("E34").Value = 0.000000028075906044477961
$ws.Range("C35").Value = 118327.9453125
$ws.Range("E35").Value = 0.000000016240404576706169
$ws.Range("C36").Value = 56474.09765625
$ws.Range("E36").Value = 0.000000015849668244527493
$ws.Range("C37").Value = 26512.328125
$ws.Range("E37").Value = 0.0000000081652329342318808
$ws.Range("C38").Value = 19543.951171875
$ws.Range("E38").Value = 0.000000010552954066156417
$ws.Range("C39").Value = 10862.0146484375
$ws.Range("E39").Value = 0.0000000061255596150999736
$ws.Range("C40").Value = 3729.3203125
$ws.Range("E40").Value = 0.0000000025664343983322624
$ws.Range("C41").Value = 1345.5067138671875
$ws.Range("E41").Value = 0.0000000017296009069411866
$ws.Range("C42").Value = 3192.2021484375
$ws.Range("E42").Value = 0.0000000029333384610197299
$ws.Range("C43").Value = 25.691064834594727
$ws.Range("E43").Value = 0.000000000033238387736611585
$ws.Range("C44").Value = 254.52424621582031
$ws.Range("E44").Value = 0.00000000031961480684294941
$ws.Range("C45").Value = 178.62748718261719
$ws.Range("E45").Value = 0.00000000049158421777661943
$ws.Range("C46").Value = 881.1624755859375
$ws.Range("E46").Value = 0.0000000024023865119460197
$ws.Range("C76").Value = 29103.9921875
$ws.Range("E76").Value = 0.000000043031484153743804
$ws.Range("C77").Value = 50032.640625
$ws.Range("E77").Value = 0.0000000060268301460553175
$ws.Range("C78").Value = 55395.890625
$ws.Range("E78").Value = 0.0000000055437436863314815
$ws.Range("C79").Value = 71883.3828125
$ws.Range("E79").Value = 0.0000000063522169746477175
$ws.Range("C80").Value = 128631.1328125
$ws.Range("E80").Value = 0.000000029041949289876356
$ws.Range("C81").Value = 223461.390625
$ws.Range("E81").Value = 0.000000052494922186951953
$ws.Range("C82").Value = 140919.203125
$ws.Range("E82").Value = 0.000000044723563519255549
$ws.Range("C83").Value = 76175.546875
$ws.Range("E83").Value = 0.000000017971478527556428
$ws.Range("C84").Value = 61683.2578125
$ws.Range("E84").Value = 0.000000010832227559376406
$ws.Range("C85").Value = 65187.88671875
$ws.Range("E85").Value = 0.0000000065915926050763574
$ws.Range("C86").Value = 33684.73828125
$ws.Range("E86").Value = 0.0000000069649495060275513
$ws.Range("C87").Value = 33559.91015625
$ws.Range("E87").Value = 0.0000000076147461669506811
$ws.Range("C88").Value = 17296.1953125
$ws.Range("E88").Value = 0.0000000068805987574194205
$ws.Range("C89").Value = 17406.576171875
$ws.Range("E89").Value = 0.0000000072320713861984132
$ws.Range("C90").Value = 8982.533203125
$ws.Range("E90").Value = 0.000000004554212118534906
$ws.Range("C91").Value = 1941.474365234375
$ws.Range("E91").Value = 0.0000000018386778766199541
$ws.Range("C92").Value = 1389.294921875
$ws.Range("E92").Value = 0.00000000094054630750406432
$ws.Range("C93").Value = 251.79917907714844
$ws.Range("E93").Value = 0.00000000024000823550807127
$ws.Range("C95").Value = 1.5880948305130005
$ws.Range("E95").Value = 0.0000000000032198824769652523
$ws.Range("C96").Value = 575.76568603515625
$ws.Range("E96").Value = 0.0000000011565025603843537
$ws.Range("C205").Value = 42157.80029296875
$ws.Range("E205").Value = 0.000000043215464984314167
$ws.Range("C206").Value = 72246.0390625
$ws.Range("E206").Value = 0.0000000060336073914868393
$ws.Range("C207").Value = 79990.453125
$ws.Range("E207").Value = 0.0000000055499778106593567
$ws.Range("C208").Value = 103798.0234375
$ws.Range("E208").Value = 0.0000000063593592614097361
$ws.Range("C209").Value = 185740.53125
$ws.Range("E209").Value = 0.000000029074605834011891
$ws.Range("C210").Value = 322673.34375
$ws.Range("E210").Value = 0.000000052553950524725224
$ws.Range("C211").Value = 203484.234375
$ws.Range("E211").Value = 0.00000004477385218137897
$ws.Range("C212").Value = 109995.8203125
$ws.Range("E212").Value = 0.000000017991686362961445
$ws.Range("C213").Value = 89069.2734375
$ws.Range("E213").Value = 0.000000010844408038224174
$ws.Range("C214").Value = 94129.8828125
$ws.Range("E214").Value = 0.0000000065990053421671746
$ws.Range("C215").Value = 48507.75
$ws.Range("E215").Value = 0.0000000069538192981610791
$ws.Range("C216").Value = 68175.4453125
$ws.Range("E216").Value = 0.000000010724822807617329
$ws.Range("C217").Value = 35136.44140625
$ws.Range("E217").Value = 0.0000000096908285840413555
$ws.Range("C218").Value = 35360.6796875
$ws.Range("E218").Value = 0.000000010185853938082801
$ws.Range("C219").Value = 18247.61328125
$ws.Range("E219").Value = 0.0000000064142811062595229
$ws.Range("C220").Value = 3944.017822265625
$ws.Range("E220").Value = 0.0000000025896456090634956
$ws.Range("C221").Value = 2822.2900390625
$ws.Range("E221").Value = 0.0000000013246919117193556
$ws.Range("C222").Value = 511.51873779296875
$ws.Range("E222").Value = 0.00000000033803435051105168
$ws.Range("C224").Value = 3.2261433601379395
$ws.Range("E224").Value = 0.0000000000045349735251776213
$ws.Range("C225").Value = 1169.64208984375
$ws.Range("E225").Value = 0.0000000016288508319917128
